# Fruta / hortaliza, semanal
# Insert 2 new weekly price rows at the top of the data block (row 26),
# pushing all existing data rows down by 2 (old row 26 -> 28, ... old row 49 -> 51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 26
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).Insert()

# ---- New row 26 ----
$ws.Cells.Item(26,1).Value = 3
$ws.Cells.Item(26,2).Value = "Femacal de La Calera"
$ws.Cells.Item(26,3).Value = "Coquimbo"
$ws.Cells.Item(26,4).Value = 45175
$ws.Cells.Item(26,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26,5).Value = 5
$ws.Cells.Item(26,6).Value = 100112043
$ws.Cells.Item(26,7).Value = "Pepino dulce"
$ws.Cells.Item(26,8).Value = "Sin especificar"
$ws.Cells.Item(26,9).Value = "Primera"
$ws.Cells.Item(26,10).Value = 68
$ws.Cells.Item(26,11).Value = 24000
$ws.Cells.Item(26,12).Value = 24000
$ws.Cells.Item(26,13).Value = 24000
$ws.Cells.Item(26,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(26,15).Value = "Provincia de Limarí"
$ws.Cells.Item(26,16).Value = 1600
$ws.Cells.Item(26,17).Value = 15
$ws.Cells.Item(26,18).Value = "Hortaliza"

# ---- New row 27 ----
$ws.Cells.Item(27,1).Value = 3
$ws.Cells.Item(27,2).Value = "Femacal de La Calera"
$ws.Cells.Item(27,3).Value = "Coquimbo"
$ws.Cells.Item(27,4).Value = 45175
$ws.Cells.Item(27,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(27,5).Value = 5
$ws.Cells.Item(27,6).Value = 100112043
$ws.Cells.Item(27,7).Value = "Pepino dulce"
$ws.Cells.Item(27,8).Value = "Sin especificar"
$ws.Cells.Item(27,9).Value = "Segunda"
$ws.Cells.Item(27,10).Value = 65
$ws.Cells.Item(27,11).Value = 19000
$ws.Cells.Item(27,12).Value = 19000
$ws.Cells.Item(27,13).Value = 19000
$ws.Cells.Item(27,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(27,15).Value = "Provincia de Limarí"
$ws.Cells.Item(27,16).Value = 1267
$ws.Cells.Item(27,17).Value = 15
$ws.Cells.Item(27,18).Value = "Hortaliza"
